$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Chris Paul / PG / San Antonio Spurs) is unchanged.

$ws.Range("A3").Value = "Russell Westbrook"
$ws.Range("C3").Value = "Denver Nuggets"

$ws.Range("A4").Value = "Brandin Podziemski"
$ws.Range("C4").Value = "Golden State Warriors"

$ws.Range("A5").Value = "Payton Pritchard"
$ws.Range("C5").Value = "Boston Celtics"

$ws.Range("A6").Value = "Deni Avdija"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Portland Trail Blazers"

$ws.Range("A7").Value = "Paolo Banchero"
$ws.Range("C7").Value = "Orlando Magic"

# Row 8 (Pascal Siakam / SF,PF,C / Indiana Pacers) is unchanged.

$ws.Range("A9").Value = "Grayson Allen"
$ws.Range("B9").Value = "PG,SG,SF"
$ws.Range("C9").Value = "Phoenix Suns"

# Row 10 (Chet Holmgren / PF,C / Oklahoma City Thunder) is unchanged.
# Row 11 (Nikola Jokic / C / Denver Nuggets) is unchanged.

$ws.Range("A12").Value = "Jerami Grant"
$ws.Range("B12").Value = "SF,PF"
$ws.Range("C12").Value = "Portland Trail Blazers"

$ws.Range("A13").Value = "Cole Anthony"
$ws.Range("B13").Value = "PG"

$ws.Range("A14").Value = "Stephon Castle"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "San Antonio Spurs"

$ws.Range("A15").Value = "Jaylen Brown"
$ws.Range("B15").Value = "SG,SF"

$ws.Range("A16").Value = "Jalen Green"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Houston Rockets"

# Row 17 (Rudy Gobert / C / Minnesota Timberwolves) is unchanged.
# Row 18 (Jalen Suggs / PG,SG / Orlando Magic) is unchanged.
# Row 19 (Jakob Poeltl / C / Toronto Raptors) is unchanged.
